$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 300
$ws.Range("I18").Value = 300
$ws.Range("K18").Value = 300
$ws.Range("M18").Value = -16

$ws.Range("H28").Value = 496.5
$ws.Range("I28").Value = 232.66667
$ws.Range("J28").Value = 760.3333
$ws.Range("K28").Value = 232.66667
$ws.Range("L28").Value = 760.3333
$ws.Range("M28").Value = 252.33333
$ws.Range("N28").Value = -1730.3333

$ws.Range("H41").Value = 280.8889
$ws.Range("I41").Value = 286.14285
$ws.Range("J41").Value = 262.5
$ws.Range("K41").Value = 286.14285
$ws.Range("L41").Value = 262.5
$ws.Range("M41").Value = 153.85715
$ws.Range("N41").Value = -1142.5

$ws.Range("H64").Value = 1799.75
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()

$ws.Range("H67").Value = 1799.75
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()

$ws.Range("H70").Value = 9844.200000000001
$ws.Range("I70").Value = 3499.875
$ws.Range("J70").Value = 35221.5
$ws.Range("K70").Value = 10499.625
$ws.Range("L70").Value = 105664.5
$ws.Range("M70").Value = -10229.625
$ws.Range("N70").Value = -106204.5

$ws.Range("H73").Value = 9844.200000000001
$ws.Range("I73").Value = 3499.875
$ws.Range("J73").Value = 35221.5
$ws.Range("K73").Value = 10499.625
$ws.Range("L73").Value = 105664.5
$ws.Range("M73").Value = -9563.625
$ws.Range("N73").Value = -107536.5

$ws.Range("H99").Value = 946.3333
$ws.Range("I99").Value = 350
$ws.Range("J99").Value = 1244.5
$ws.Range("K99").Value = 1050
$ws.Range("L99").Value = 3733.5
$ws.Range("M99").Value = 448
$ws.Range("N99").Value = -6729.5

$ws.Range("H118").Value = 190
$ws.Range("I118").Value = 190
$ws.Range("K118").Value = 570
$ws.Range("M118").Value = 1087

$ws.Range("H127").Value = 848.5
$ws.Range("I127").Value = 848.5
$ws.Range("K127").Value = 2545.5
$ws.Range("M127").Value = 2414.5

$ws.Range("H132").Value = 1105.125
$ws.Range("I132").Value = 1125.0667
$ws.Range("J132").Value = 806
$ws.Range("K132").Value = 3375.2001
$ws.Range("L132").Value = 2418
$ws.Range("M132").Value = -845.2001
$ws.Range("N132").Value = -7478

$ws.Range("H137").Value = 1263.8182
$ws.Range("I137").Value = 1300.2858
$ws.Range("K137").Value = 3900.8574
$ws.Range("M137").Value = -1350.8574

$ws.Range("H138").Value = 2071.3845
$ws.Range("I138").Value = 1469
$ws.Range("J138").Value = 5384.5
$ws.Range("K138").Value = 4407
$ws.Range("L138").Value = 16153.5
$ws.Range("M138").Value = 733
$ws.Range("N138").Value = -26433.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 212.5
$ws.Range("I4").Value = 324.5
$ws.Range("J4").Value = 100.5
$ws.Range("K4").Value = 324.5
$ws.Range("L4").Value = 100.5
$ws.Range("M4").Value = -208.5
$ws.Range("N4").Value = -332.5

$ws.Range("H5").Value = 393.5
$ws.Range("I5").Value = 274.66666
$ws.Range("K5").Value = 274.66666
$ws.Range("M5").Value = -162.66666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 393.5
$ws.Range("I4").Value = 274.66666
$ws.Range("K4").Value = 274.66666
$ws.Range("M4").Value = -159.66666

$ws.Range("H20").Value = 2125
$ws.Range("I20").Value = 1800
$ws.Range("K20").Value = 1800
$ws.Range("M20").Value = -1553

$ws.Range("H134").Value = 2128.4285
$ws.Range("I134").Value = 2128.4285
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 6385.2855
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -3850.2855
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 86.09677000000001
$ws.Range("I7").Value = 81.72221999999999
$ws.Range("J7").Value = 92.15385000000001
$ws.Range("K7").Value = 81.72221999999999
$ws.Range("L7").Value = 92.15385000000001
$ws.Range("M7").Value = 31.27778000000001
$ws.Range("N7").Value = -318.15385

$ws.Range("H22").Value = 99
$ws.Range("I22").Value = 99
$ws.Range("K22").Value = 99
$ws.Range("M22").Value = 251

$ws.Range("H31").Value = 1452.3846
$ws.Range("I31").Value = 1452.3846
$ws.Range("K31").Value = 1452.3846
$ws.Range("M31").Value = -1157.3846

$ws.Range("H34").Value = 1452.3846
$ws.Range("I34").Value = 1452.3846
$ws.Range("K34").Value = 1452.3846
$ws.Range("M34").Value = -1250.3846

$ws.Range("H45").Value = 38000
$ws.Range("J45").Value = 38000
$ws.Range("L45").Value = 38000
$ws.Range("N45").Value = -39186

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 66.90909000000001
$ws.Range("I2").Value = 29.142857
$ws.Range("J2").Value = 133
$ws.Range("K2").Value = 174.857142
$ws.Range("L2").Value = 798
$ws.Range("M2").Value = -61.85714200000001
$ws.Range("N2").Value = -1024

$ws.Range("H59").Value = 571.4286
$ws.Range("I59").Value = 500
$ws.Range("K59").Value = 1500
$ws.Range("M59").Value = -960

$ws.Range("H64").Value = 4220
$ws.Range("J64").Value = 4220
$ws.Range("L64").Value = 12660
$ws.Range("N64").Value = -13200

$ws.Range("H67").Value = 4220
$ws.Range("J67").Value = 4220
$ws.Range("L67").Value = 12660
$ws.Range("N67").Value = -14532

$ws.Range("H104").Value = 500
$ws.Range("I104").Value = 500
$ws.Range("K104").Value = 1500
$ws.Range("M104").Value = 1121

$ws.Range("H107").Value = 800
$ws.Range("I107").Value = 400
$ws.Range("K107").Value = 1200
$ws.Range("M107").Value = 720

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 1632
$ws.Range("I5").Value = 1632
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 1632
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -1520
$ws.Range("N5").ClearContents()

$ws.Range("H70").Value = 3699
$ws.Range("I70").Value = 3699
$ws.Range("K70").Value = 3699
$ws.Range("M70").Value = -3429

$ws.Range("H73").Value = 3699
$ws.Range("I73").Value = 3699
$ws.Range("K73").Value = 3699
$ws.Range("M73").Value = -2763

$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

$ws.Range("H132").Value = 4924.75
$ws.Range("I132").Value = 3849.5
$ws.Range("J132").Value = 6000
$ws.Range("K132").Value = 11548.5
$ws.Range("L132").Value = 18000
$ws.Range("M132").Value = -9018.5
$ws.Range("N132").Value = -23060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()

$ws.Range("H33").Value = 5000
$ws.Range("I33").Value = 5000
$ws.Range("K33").Value = 5000
$ws.Range("M33").Value = -4710

$ws.Range("H68").Value = 1661.3334
$ws.Range("I68").Value = 1661.3334
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 1661.3334
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -912.3334
$ws.Range("N68").ClearContents()

$ws.Range("H69").Value = 52000
$ws.Range("J69").Value = 52000
$ws.Range("L69").Value = 52000
$ws.Range("N69").Value = -53622

$ws.Range("H71").Value = 1661.3334
$ws.Range("I71").Value = 1661.3334
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 8306.666999999999
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -4562.666999999999
$ws.Range("N71").ClearContents()

$ws.Range("H72").Value = 52000
$ws.Range("J72").Value = 52000
$ws.Range("L72").Value = 156000
$ws.Range("N72").Value = -164112

$ws.Range("H104").Value = 42750
$ws.Range("J104").Value = 42750
$ws.Range("L104").Value = 42750
$ws.Range("N104").Value = -49738

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 251100
$ws.Range("I2").Value = 334300
$ws.Range("K2").Value = 334300
$ws.Range("M2").Value = -334188

$ws.Range("H107").Value = 536.82355
$ws.Range("I107").Value = 244.08333
$ws.Range("K107").Value = 732.24999
$ws.Range("M107").Value = 1187.75001

$ws.Range("H132").Value = 2269.6453
$ws.Range("I132").Value = 1728.1111
$ws.Range("J132").Value = 5925
$ws.Range("K132").Value = 5184.3333
$ws.Range("L132").Value = 17775
$ws.Range("M132").Value = -2654.3333
$ws.Range("N132").Value = -22835
